$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = "60/66"
$ws.Range("H3").Value = "0/66"
$ws.Range("H4").Value = "66/66"
$ws.Range("H5").Value = "51/66"
$ws.Range("H6").Value = "0/66"
$ws.Range("H7").Value = "59/66"
$ws.Range("H8").Value = "60/66"
$ws.Range("H9").Value = "66/66"
$ws.Range("H10").Value = "40/66"
$ws.Range("H11").Value = "27/66"
$ws.Range("H12").Value = "53/66"
$ws.Range("H13").Value = "0/66"
$ws.Range("H14").Value = "47/66"
$ws.Range("H15").Value = "25/66"
$ws.Range("H16").Value = "51/66"
$ws.Range("H17").Value = "58/66"
$ws.Range("H18").Value = "58/66"
$ws.Range("H19").Value = "0/66"
$ws.Range("H20").Value = "28/66"
$ws.Range("H21").Value = "31/66"
$ws.Range("H22").Value = "45/66"
$ws.Range("H23").Value = "38/66"
$ws.Range("H24").Value = "45/66"
$ws.Range("H25").Value = "0/66"
$ws.Range("H26").Value = "66/66"
$ws.Range("H27").Value = "62/66"
$ws.Range("H28").Value = "55/66"
$ws.Range("H29").Value = "59/66"
$ws.Range("H30").Value = "0/66"
$ws.Range("H31").Value = "0/66"
$ws.Range("H32").Value = "57/66"
$ws.Range("H33").Value = "60/66"
$ws.Range("H34").Value = "52/66"
$ws.Range("H35").Value = "30/66"
$ws.Range("H36").Value = "0/66"
$ws.Range("H37").Value = "45/66"
$ws.Range("H38").Value = "57/66"
$ws.Range("H39").Value = "59/66"
$ws.Range("H40").Value = "27/66"
$ws.Range("H41").Value = "17/66"
$ws.Range("H42").Value = "47/66"
$ws.Range("H43").Value = "54/66"
$ws.Range("H44").Value = "60/66"
$ws.Range("H45").Value = "55/66"
$ws.Range("H46").Value = "42/66"
$ws.Range("H47").Value = "1/66"
$ws.Range("H48").Value = "59/66"
$ws.Range("H49").Value = "52/66"
$ws.Range("H50").Value = "54/66"
$ws.Range("H51").Value = "43/66"
$ws.Range("H52").Value = "0/66"
$ws.Range("H53").Value = "0/66"
$ws.Range("H54").Value = "0/66"
$ws.Range("H55").Value = "0/66"
$ws.Range("H56").Value = "0/66"
$ws.Range("H112").Value = "47/64"
$ws.Range("H113").Value = "57/64"
$ws.Range("H114").Value = "60/64"
$ws.Range("H115").Value = "55/64"
$ws.Range("H116").Value = "0/64"
$ws.Range("H117").Value = "62/64"
$ws.Range("H118").Value = "59/64"
$ws.Range("H119").Value = "52/64"
$ws.Range("H120").Value = "57/64"
$ws.Range("H121").Value = "35/64"
$ws.Range("H122").Value = "56/64"
$ws.Range("H123").Value = "59/64"
$ws.Range("H124").Value = "0/64"
$ws.Range("H125").Value = "51/64"
$ws.Range("H126").Value = "0/64"
$ws.Range("H127").Value = "0/64"
$ws.Range("H128").Value = "0/64"
$ws.Range("H129").Value = "60/64"
$ws.Range("H130").Value = "42/64"
$ws.Range("H131").Value = "40/64"
$ws.Range("H132").Value = "0/64"
$ws.Range("H133").Value = "0/64"
$ws.Range("H134").Value = "0/64"
$ws.Range("H135").Value = "0/64"
$ws.Range("H136").Value = "0/64"
$ws.Range("H137").Value = "61/64"
$ws.Range("H138").Value = "60/64"
$ws.Range("H139").Value = "63/64"
$ws.Range("H140").Value = "53/64"
$ws.Range("H141").Value = "49/64"
$ws.Range("H142").Value = "53/64"
$ws.Range("H143").Value = "32/64"
$ws.Range("H144").Value = "25/64"
$ws.Range("H145").Value = "0/64"
$ws.Range("H146").Value = "0/64"
$ws.Range("H147").Value = "0/64"
$ws.Range("H148").Value = "0/64"
$ws.Range("H149").Value = "0/64"
$ws.Range("H150").Value = "0/64"
$ws.Range("H151").Value = "0/64"
$ws.Range("H152").Value = "56/64"
$ws.Range("H153").Value = "60/64"
$ws.Range("H154").Value = "60/64"
$ws.Range("H155").Value = "0/64"
$ws.Range("H156").Value = "0/64"
$ws.Range("H157").Value = "0/64"
$ws.Range("H158").Value = "55/64"
$ws.Range("H159").Value = "55/64"
$ws.Range("H160").Value = "45/64"
$ws.Range("H161").Value = "0/64"
$ws.Range("H162").Value = "54/64"
$ws.Range("H163").Value = "59/64"
$ws.Range("H164").Value = "32/64"
$ws.Range("H165").Value = "60/64"
$ws.Range("H166").Value = "52/64"
$ws.Range("H167").Value = "0/62"
$ws.Range("H168").Value = "0/62"
$ws.Range("H169").Value = "0/62"
$ws.Range("H170").Value = "0/62"
$ws.Range("H171").Value = "0/62"
$ws.Range("H172").Value = "55/62"
$ws.Range("H173").Value = "54/62"
$ws.Range("H174").Value = "0/62"
$ws.Range("H175").Value = "0/62"
$ws.Range("H176").Value = "0/62"
$ws.Range("H177").Value = "1/62"
$ws.Range("H178").Value = "54/62"
$ws.Range("H179").Value = "51/62"
$ws.Range("H180").Value = "33/62"
$ws.Range("H181").Value = "60/62"
$ws.Range("H182").Value = "52/62"
$ws.Range("H183").Value = "47/62"
$ws.Range("H184").Value = "56/62"
$ws.Range("H185").Value = "60/62"
$ws.Range("H186").Value = "33/62"
$ws.Range("H187").Value = "43/62"
$ws.Range("H188").Value = "0/62"
$ws.Range("H189").Value = "46/62"
$ws.Range("H190").Value = "0/62"
$ws.Range("H191").Value = "0/62"
$ws.Range("H192").Value = "62/62"
$ws.Range("H193").Value = "45/62"
$ws.Range("H194").Value = "40/62"
$ws.Range("H195").Value = "57/62"
$ws.Range("H196").Value = "62/62"
$ws.Range("H197").Value = "53/62"
$ws.Range("H198").Value = "57/62"
$ws.Range("H199").Value = "53/62"
$ws.Range("H200").Value = "0/62"
$ws.Range("H201").Value = "56/62"
$ws.Range("H202").Value = "48/62"
$ws.Range("H203").Value = "56/62"
$ws.Range("H204").Value = "31/62"
$ws.Range("H205").Value = "42/62"
$ws.Range("H206").Value = "48/62"
$ws.Range("H207").Value = "49/62"
$ws.Range("H208").Value = "53/62"
$ws.Range("H209").Value = "44/62"
$ws.Range("H210").Value = "55/62"
$ws.Range("H211").Value = "40/62"
$ws.Range("H212").Value = "57/62"
$ws.Range("H213").Value = "10/62"
$ws.Range("H214").Value = "58/62"
$ws.Range("H215").Value = "58/62"
$ws.Range("H216").Value = "51/62"
$ws.Range("H217").Value = "53/62"
$ws.Range("H218").Value = "55/62"
$ws.Range("H219").Value = "54/62"
$ws.Range("H220").Value = "49/62"
$ws.Range("H221").Value = "32/62"
$ws.Range("H222").Value = "43/60"
$ws.Range("H223").Value = "31/60"
$ws.Range("H224").Value = "50/60"
$ws.Range("H225").Value = "38/60"
$ws.Range("H226").Value = "0/60"
$ws.Range("H227").Value = "0/60"
$ws.Range("H228").Value = "55/60"
$ws.Range("H229").Value = "0/60"
$ws.Range("H230").Value = "50/60"
$ws.Range("H231").Value = "35/60"
$ws.Range("H232").Value = "40/60"
$ws.Range("H233").Value = "49/60"
$ws.Range("H234").Value = "54/60"
$ws.Range("H235").Value = "0/60"
$ws.Range("H236").Value = "0/60"
$ws.Range("H237").Value = "0/60"
$ws.Range("H238").Value = "0/60"
$ws.Range("H239").Value = "0/60"
$ws.Range("H240").Value = "0/60"
$ws.Range("H241").Value = "0/60"
$ws.Range("H242").Value = "37/60"
$ws.Range("H243").Value = "38/60"
$ws.Range("H244").Value = "31/60"
$ws.Range("H245").Value = "44/60"
$ws.Range("H246").Value = "0/60"
$ws.Range("H247").Value = "41/60"
$ws.Range("H248").Value = "51/60"
$ws.Range("H249").Value = "50/60"
$ws.Range("H250").Value = "43/60"
$ws.Range("H251").Value = "0/60"
$ws.Range("H252").Value = "22/60"
$ws.Range("H253").Value = "49/60"
$ws.Range("H254").Value = "29/60"
$ws.Range("H255").Value = "56/60"
$ws.Range("H256").Value = "23/60"
$ws.Range("H257").Value = "23/60"
$ws.Range("H258").Value = "51/60"
$ws.Range("H259").Value = "0/60"
$ws.Range("H260").Value = "43/60"
$ws.Range("H261").Value = "25/60"
$ws.Range("H262").Value = "35/60"
$ws.Range("H263").Value = "52/60"
$ws.Range("H264").Value = "42/60"
$ws.Range("H265").Value = "0/60"
$ws.Range("H266").Value = "54/60"
$ws.Range("H267").Value = "45/60"
$ws.Range("H268").Value = "58/60"
$ws.Range("H269").Value = "54/60"
$ws.Range("H270").Value = "22/60"
$ws.Range("H271").Value = "38/60"
$ws.Range("H272").Value = "48/60"
$ws.Range("H273").Value = "55/60"
$ws.Range("H274").Value = "42/60"
$ws.Range("H275").Value = "44/60"
$ws.Range("H276").Value = "40/60"
$ws.Range("L4").Value = 314
$ws.Range("L10").NumberFormat = "@"
$ws.Range("L10").Value = "75.8%"
$ws.Range("L10").NumberFormat = "General"
$ws.Range("M15").Value = 66
$ws.Range("S15").NumberFormat = "@"
$ws.Range("S15").Value = "72.9%"
$ws.Range("S15").NumberFormat = "General"
$ws.Range("M17").Value = 64
$ws.Range("S17").NumberFormat = "@"
$ws.Range("S17").Value = "81.6%"
$ws.Range("S17").NumberFormat = "General"
$ws.Range("M18").Value = 62
$ws.Range("S18").NumberFormat = "@"
$ws.Range("S18").Value = "77.8%"
$ws.Range("S18").NumberFormat = "General"
$ws.Range("M19").Value = 60
$ws.Range("S19").NumberFormat = "@"
$ws.Range("S19").Value = "70.3%"
$ws.Range("S19").NumberFormat = "General"
$ws.Range("G125").Value = "abdallahashraf2023@gmail.com, ahmedali78112@gmail.com"
$ws.Range("G222").Value = "fatma_shoukry@hotmail.com, drmohamedramadan50@gmail.com"
